$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 47438
$ws.Range("C10").Value = 'SIG-3w Lilliput LED Torch &amp; Table Lamp'
$ws.Range("D10").Value = 401.81
$ws.Range("E10").Value = 480.05
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 803.62
$ws.Range("B11").Value = 59408
$ws.Range("C11").Value = 'SIG-3W Lilliput LED Torch &amp; Table Lamp'
$ws.Range("D11").Value = 388.17
$ws.Range("E11").Value = 463.78
$ws.Range("F11").Value = 13
$ws.Range("G11").Value = 5046.21
$ws.Range("F26").Value = 89
$ws.Range("G26").Value = 4102.9
$ws.Range("B46").Value = 29640.89
$ws.Range("F83").Value = 73
$ws.Range("G83").Value = 2442.58
$ws.Range("B85").Value = 168487.01
$ws.Range("F95").Value = 15
$ws.Range("G95").Value = 3768.75
$ws.Range("F97").Value = 10
$ws.Range("G97").Value = 3771.9
$ws.Range("F100").Value = 17
$ws.Range("G100").Value = 2023.34
$ws.Range("B103").Value = 37646.36
$ws.Range("F218").Value = 30
$ws.Range("G218").Value = 2347.8
$ws.Range("F219").Value = 209
$ws.Range("G219").Value = 6485.27
$ws.Range("F224").Value = 52
$ws.Range("G224").Value = 4226.56
$ws.Range("B228").Value = 37006.93
$ws.Range("F266").Value = 42
$ws.Range("G266").Value = 1760.22
$ws.Range("B267").Value = 23365.69
$ws.Range("F294").Value = 18
$ws.Range("G294").Value = 14556.24
$ws.Range("B298").Value = 116854.81
$ws.Range("F324").Value = 47
$ws.Range("G324").Value = 3681.51
$ws.Range("F332").Value = 77
$ws.Range("G332").Value = 2564.1
$ws.Range("F343").Value = 23
$ws.Range("G343").Value = 1078.01
$ws.Range("F345").Value = 42
$ws.Range("G345").Value = 3670.38
$ws.Range("B349").Value = 145049.9
$ws.Range("F361").Value = 67
$ws.Range("G361").Value = 4938.57
$ws.Range("F372").Value = 34
$ws.Range("G372").Value = 1619.76
$ws.Range("F382").Value = 40
$ws.Range("G382").Value = 4392.8
$ws.Range("F392").Value = 7
$ws.Range("G392").Value = 600.04
$ws.Range("F398").Value = 5
$ws.Range("G398").Value = 773.3
$ws.Range("F402").Value = 34
$ws.Range("G402").Value = 1852.32
$ws.Range("F409").Value = 232
$ws.Range("G409").Value = 39748.56
$ws.Range("F415").Value = 6
$ws.Range("G415").Value = 793.38
$ws.Range("F418").Value = 18
$ws.Range("G418").Value = 1070.46
$ws.Range("F419").Value = 305
$ws.Range("G419").Value = 12559.9
$ws.Range("F420").Value = 90
$ws.Range("G420").Value = 8336.700000000001
$ws.Range("F422").Value = 28
$ws.Range("G422").Value = 4022.48
$ws.Range("B423").Value = 182686.19
$ws.Range("F425").Value = 20
$ws.Range("G425").Value = 3671.6
$ws.Range("F431").Value = 19
$ws.Range("G431").Value = 4285.83
$ws.Range("B437").Value = 27142.93
$ws.Range("F449").Value = 35
$ws.Range("G449").Value = 1685.25
$ws.Range("B455").Value = 24779.57
$ws.Range("F469").Value = 39
$ws.Range("G469").Value = 3626.61
$ws.Range("F479").Value = 25
$ws.Range("G479").Value = 3616
$ws.Range("B481").Value = 47117.75
$ws.Range("F487").Value = 7
$ws.Range("G487").Value = 387.31
$ws.Range("F496").Value = 108
$ws.Range("G496").Value = 16229.16
$ws.Range("B497").Value = 41432.75
$ws.Range("B518").Value = 58047
$ws.Range("D518").Value = 105.54
$ws.Range("E518").Value = 126.1
$ws.Range("F518").Value = 71
$ws.Range("G518").Value = 7493.34
$ws.Range("B519").Value = 47097
$ws.Range("D519").Value = 112.28
$ws.Range("E519").Value = 134.16
$ws.Range("F519").Value = 18
$ws.Range("G519").Value = 2021.04
$ws.Range("F520").Value = 141
$ws.Range("G520").Value = 1386.03
$ws.Range("F521").Value = 178
$ws.Range("G521").Value = 4862.96
$ws.Range("F524").Value = 198
$ws.Range("G524").Value = 5888.52
$ws.Range("F526").Value = 830
$ws.Range("G526").Value = 80178
$ws.Range("F529").Value = 165
$ws.Range("G529").Value = 4436.85
$ws.Range("B532").Value = 160833.86
$ws.Range("F561").Value = 947
$ws.Range("G561").Value = 12216.3
$ws.Range("F566").Value = 180
$ws.Range("G566").Value = 3457.8
$ws.Range("B567").Value = 54873.21
$ws.Range("F591").Value = 348
$ws.Range("G591").Value = 12865.56
$ws.Range("B610").Value = 62837.4
$ws.Range("F614").Value = 17
$ws.Range("G614").Value = 3872.26
$ws.Range("F621").Value = 258
$ws.Range("G621").Value = 15665.76
$ws.Range("F625").Value = 51
$ws.Range("G625").Value = 3280.32
$ws.Range("B638").Value = 155033.94
$ws.Range("F671").Value = 174
$ws.Range("G671").Value = 10770.6
$ws.Range("F672").Value = 79
$ws.Range("G672").Value = 21021.11
$ws.Range("F674").Value = 99
$ws.Range("G674").Value = 5165.82
$ws.Range("F684").Value = 42
$ws.Range("G684").Value = 7279.44
$ws.Range("B688").Value = 93200.55
$ws.Range("F713").Value = 37
$ws.Range("G713").Value = 4830.35
$ws.Range("F715").Value = 60
$ws.Range("G715").Value = 7833
$ws.Range("F717").Value = 66
$ws.Range("G717").Value = 1795.2
$ws.Range("F718").Value = 163
$ws.Range("G718").Value = 4433.6
$ws.Range("F719").Value = 134
$ws.Range("G719").Value = 3644.8
$ws.Range("B720").Value = 40410.41
$ws.Range("F747").Value = 3
$ws.Range("G747").Value = 4602.9
$ws.Range("B773").Value = 158768.65
$ws.Range("F812").Value = 62
$ws.Range("G812").Value = 4140.98
$ws.Range("B815").Value = 39535.48
$ws.Range("F817").Value = 14
$ws.Range("G817").Value = 1141.84
$ws.Range("F821").Value = 150
$ws.Range("G821").Value = 19965
$ws.Range("F829").Value = 112
$ws.Range("G829").Value = 6017.76
$ws.Range("F833").Value = 139
$ws.Range("G833").Value = 15405.37
$ws.Range("B837").Value = 206612.27
$ws.Range("F844").Value = 17
$ws.Range("G844").Value = 430.61
$ws.Range("F861").Value = 317
$ws.Range("G861").Value = 11675.11
$ws.Range("F865").Value = 118
$ws.Range("G865").Value = 5887.02
$ws.Range("B867").Value = 222933.89
$ws.Range("B923").Value = 2948990.34
$ws.Range("B924").Value = 2948990.34

Write-Host "Applied 170 cell updates"
